# File Area: Fix support for .xlsx files
#
# The original sheet only had data in columns B:D (rows 2-5). The fix adds
# column A to the sheet: a new header-less row 1 (A1:D1) formatted like the
# rest of the table, a blank formatted cell in column A for rows 2-4, and a
# numeric value (5) in A5 - all using the same "label" style already used by
# cells such as C3/C4/C5 (centered, Times New Roman).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting already used for the other label cells (e.g. C3) and
# apply it (format only) to the new column A cells and the new row 1, so the
# existing style index is reused instead of creating a new one.
$ws.Range("C3").Copy()
[void]$ws.Range("A1:D1").PasteSpecial(-4122)
[void]$ws.Range("A2").PasteSpecial(-4122)
[void]$ws.Range("A3").PasteSpecial(-4122)
[void]$ws.Range("A4").PasteSpecial(-4122)
[void]$ws.Range("A5").PasteSpecial(-4122)

# A5 gets an actual numeric value, the rest of the new cells stay empty.
$ws.Range("A5").Value = 5

# Update the selected cell to B5, matching the saved selection state.
[void]$ws.Range("B5").Select()
